$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New speaker_variant values (column C) and derived id values (column B),
# re-exported without is_prefered ("x") markers.
$rows = @(
    @{ Row = 2;  Id = "#petiet";        Name = "Petiet" },
    @{ Row = 3;  Id = "#baggalyn";      Name = "Baggalyn" },
    @{ Row = 4;  Id = "#mansh";         Name = "Mansh" },
    @{ Row = 5;  Id = "#pronk";         Name = "Pronk" },
    @{ Row = 6;  Id = "#manshart";      Name = "Manshart" },
    @{ Row = 7;  Id = "#narticoforus";  Name = "Narticoforus" },
    @{ Row = 8;  Id = "#bartrand";      Name = "Bartrand" },
    @{ Row = 9;  Id = "#mamshart";      Name = "Mamshart" },
    @{ Row = 10; Id = "#alet";          Name = "Alet" },
    @{ Row = 11; Id = "#manshert";      Name = "Manshert" },
    @{ Row = 12; Id = "#florentyn";     Name = "Florentyn" },
    @{ Row = 13; Id = "#helena";        Name = "Helena" },
    @{ Row = 14; Id = "#anna";          Name = "Anna" },
    @{ Row = 15; Id = "#boudewyn";      Name = "Boudewyn" },
    @{ Row = 16; Id = "#tryn";          Name = "Tryn" },
    @{ Row = 17; Id = "#broer";         Name = "Broer" },
    @{ Row = 18; Id = "#ferdinand";     Name = "Ferdinand" },
    @{ Row = 19; Id = "#marticoforus";  Name = "Marticoforus" },
    @{ Row = 20; Id = "#bagotyn";       Name = "Bagotyn" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Id
    $ws.Cells.Item($r.Row, 3).Value = $r.Name
    $ws.Cells.Item($r.Row, 4).Value = ""
}
